# Applies the "Add files via upload" commit:
#   - adds a new "cdmContracts" lookup sheet (+ its Table4 list object) just
#     before the "XtraButtons" sheet
#   - replaces the old "Residual Risk Owner 1/2/3" list on cdmResidualRiskOwners
#     with the new two-item "HS2 ..." list
#   - adds a "Sync HS2 Hazards" / "synccsv" row to XtraButtons

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. cdmResidualRiskOwners: shrink the Residual Risk Owner list from three
#    rows to two, with new content.
# ---------------------------------------------------------------------
$owners = $wb.Worksheets.Item("cdmResidualRiskOwners")
$owners.Range("A2").Value = "HS2 Infrastructure Management SME"
$owners.Range("A3").Value = "HS2 Rail Systems Interface Engineer"
$owners.Rows(4).Delete()
$owners.Range("A4").Select()

# ---------------------------------------------------------------------
# 2. New "cdmContracts" sheet, inserted right before "XtraButtons".
# ---------------------------------------------------------------------
$xtraButtonsSheet = $wb.Worksheets.Item("XtraButtons")
$contracts = $wb.Worksheets.Add($xtraButtonsSheet)
$contracts.Name = "cdmContracts"
$contracts.Tab.Color = 255

$contractTitles = @(
  "Title",
  "HRS06 - NICC & Washwood Heath Depot",
  "HRS07 - Precast Slab Track System",
  "HRS11 - Cross Passageway Doors",
  "HRS12 - OCS",
  "HRS13 - Track Urban Phase 1",
  "HRS14 - Track Open Route Central Phase 1 (incl. Calvert IMD)",
  "HRS15 - Track Open Route North Phase 1",
  "HRS16 - Track Phase 2a (incl. Stone IMB-R)",
  "HRS17 - Rail - Supply",
  "HRS18 - S&C - Supply",
  "HRS19 - Maintenance HV & Traction Power",
  "HRS20 - Operational Telecommunications",
  "HRS21 - 3rd Party Telecommunications",
  "HRS22 - Engineering Management System",
  "HRS23 - CCS and TM",
  "HRS24 - Network Rail - Communications Contract Placeholder"
)

for ($i = 0; $i -lt $contractTitles.Length; $i++) {
  $contracts.Cells.Item($i + 1, 1).Value = $contractTitles[$i]
}

$contracts.Columns(1).ColumnWidth = 52.63

$contractsTable = $contracts.ListObjects.Add(1, $contracts.Range("A1:A17"), 0, 1)
$contractsTable.Name = "Table4"

$contracts.Range("G14").Select()

# ---------------------------------------------------------------------
# 3. XtraButtons: add the "Sync HS2 Hazards" / "synccsv" button row.
# ---------------------------------------------------------------------
$buttons = $wb.Worksheets.Item("XtraButtons")
$buttons.Range("A3").Value = "Sync HS2 Hazards"
$buttons.Range("B3").Value = "synccsv"
$buttons.Range("A3").Select()
